$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2595.3333
$ws.Range("I6").Value = 103.90909
$ws.Range("J6").Value = 30001
$ws.Range("K6").Value = 311.72727
$ws.Range("L6").Value = 90003
$ws.Range("M6").Value = -199.72727
$ws.Range("N6").Value = -90227
$ws.Range("H8").Value = 61.2
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = $null
$ws.Range("H17").Value = 1011951.75
$ws.Range("J17").Value = 1011951.75
$ws.Range("L17").Value = 3035855.25
$ws.Range("N17").Value = -3036191.25
$ws.Range("H19").Value = 1811.8334
$ws.Range("I19").Value = 1420.7273
$ws.Range("J19").Value = 2142.7693
$ws.Range("K19").Value = 1420.7273
$ws.Range("L19").Value = 2142.7693
$ws.Range("M19").Value = -1245.7273
$ws.Range("N19").Value = -2492.7693
$ws.Range("H70").Value = 1647.1177
$ws.Range("J70").Value = 1647.1177
$ws.Range("L70").Value = 4941.3531
$ws.Range("N70").Value = -5481.3531
$ws.Range("H73").Value = 1647.1177
$ws.Range("J73").Value = 1647.1177
$ws.Range("L73").Value = 4941.3531
$ws.Range("N73").Value = -6813.3531
$ws.Range("H76").Value = 11805.375
$ws.Range("I76").Value = 14088.6
$ws.Range("K76").Value = 14088.6
$ws.Range("M76").Value = -13773.6
$ws.Range("H79").Value = 11805.375
$ws.Range("I79").Value = 14088.6
$ws.Range("K79").Value = 14088.6
$ws.Range("M79").Value = -12996.6
$ws.Range("H80").Value = 1578.0322
$ws.Range("I80").Value = 354.83334
$ws.Range("J80").Value = 2350.5789
$ws.Range("K80").Value = 1064.50002
$ws.Range("L80").Value = 7051.736699999999
$ws.Range("M80").Value = -66.50001999999995
$ws.Range("N80").Value = -9047.736699999999
$ws.Range("H83").Value = 1578.0322
$ws.Range("I83").Value = 354.83334
$ws.Range("J83").Value = 2350.5789
$ws.Range("K83").Value = 3193.50006
$ws.Range("L83").Value = 21155.2101
$ws.Range("M83").Value = 1798.49994
$ws.Range("N83").Value = -31139.2101
$ws.Range("H101").Value = 1056.6666
$ws.Range("J101").Value = 1692.5
$ws.Range("L101").Value = 5077.5
$ws.Range("N101").Value = -8321.5
$ws.Range("H113").Value = 3565.2222
$ws.Range("I113").Value = 2479.818
$ws.Range("J113").Value = 5270.857
$ws.Range("K113").Value = 2479.818
$ws.Range("L113").Value = 5270.857
$ws.Range("M113").Value = 774.1819999999998
$ws.Range("N113").Value = -11778.857
$ws.Range("H115").Value = 384.625
$ws.Range("I115").Value = 368.14285
$ws.Range("K115").Value = 1104.42855
$ws.Range("M115").Value = 462.5714499999999
$ws.Range("H132").Value = 3645.4482
$ws.Range("I132").Value = 3366.652
$ws.Range("K132").Value = 10099.956
$ws.Range("M132").Value = -7569.956
$ws.Range("H135").Value = 41668680
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 41668680
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 375018120
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = -375023190
$ws.Range("H137").Value = 2272.0356
$ws.Range("I137").Value = 2259.842
$ws.Range("J137").Value = 2297.7778
$ws.Range("K137").Value = 6779.526
$ws.Range("L137").Value = 6893.3334
$ws.Range("M137").Value = -4229.526
$ws.Range("N137").Value = -11993.3334
$ws.Range("H138").Value = 7411256.5
$ws.Range("I138").Value = 994.8570999999999
$ws.Range("J138").Value = 13895235
$ws.Range("K138").Value = 2984.5713
$ws.Range("L138").Value = 41685705
$ws.Range("M138").Value = 2155.4287
$ws.Range("N138").Value = -41695985
$ws.Range("H141").Value = 2235.611
$ws.Range("I141").Value = 1702.5625
$ws.Range("J141").Value = 6500
$ws.Range("K141").Value = 5107.6875
$ws.Range("L141").Value = 19500
$ws.Range("M141").Value = 72.3125
$ws.Range("N141").Value = -29860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29423638
$ws.Range("I32").Value = 45461380
$ws.Range("J32").Value = 21112.084
$ws.Range("K32").Value = 45461380
$ws.Range("L32").Value = 21112.084
$ws.Range("M32").Value = -45461093
$ws.Range("N32").Value = -21686.084
$ws.Range("H45").Value = 2584.2856
$ws.Range("J45").Value = 3347.3333
$ws.Range("L45").Value = 3347.3333
$ws.Range("N45").Value = -4101.3333
$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 5000
$ws.Range("M60").Value = -4267
$ws.Range("H61").Value = 71434320
$ws.Range("I61").Value = 111114450
$ws.Range("J61").Value = 10085.6
$ws.Range("K61").Value = 111114450
$ws.Range("L61").Value = 10085.6
$ws.Range("M61").Value = -111114238
$ws.Range("N61").Value = -10509.6
$ws.Range("H74").Value = 200004030
$ws.Range("I74").Value = 500002500
$ws.Range("J74").Value = 5066
$ws.Range("K74").Value = 500002500
$ws.Range("L74").Value = 5066
$ws.Range("M74").Value = -500001626
$ws.Range("N74").Value = -6814
$ws.Range("H77").Value = 200004030
$ws.Range("I77").Value = 500002500
$ws.Range("J77").Value = 5066
$ws.Range("K77").Value = 2500012500
$ws.Range("L77").Value = 25330
$ws.Range("M77").Value = -2500008132
$ws.Range("N77").Value = -34066
$ws.Range("H102").Value = 2312.1052
$ws.Range("I102").Value = 1621.375
$ws.Range("K102").Value = 1621.375
$ws.Range("M102").Value = 0.625
$ws.Range("H132").Value = 52634412
$ws.Range("I132").Value = 2988
$ws.Range("K132").Value = 8964
$ws.Range("M132").Value = -6434
$ws.Range("H134").Value = 59832.25
$ws.Range("J134").Value = 59832.25
$ws.Range("L134").Value = 59832.25
$ws.Range("N134").Value = -69972.25
$ws.Range("H136").Value = 71434320
$ws.Range("I136").Value = 111114450
$ws.Range("J136").Value = 10085.6
$ws.Range("K136").Value = 333343350
$ws.Range("L136").Value = 30256.8
$ws.Range("M136").Value = -333340800
$ws.Range("N136").Value = -35356.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1161.8276
$ws.Range("I94").Value = 892.2857
$ws.Range("K94").Value = 892.2857
$ws.Range("M94").Value = -441.2857
$ws.Range("H99").Value = 6335.1665
$ws.Range("J99").Value = 7377.75
$ws.Range("L99").Value = 7377.75
$ws.Range("N99").Value = -10373.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 311.1
$ws.Range("I7").Value = 232.7
$ws.Range("K7").Value = 232.7
$ws.Range("M7").Value = -119.7
$ws.Range("H52").Value = 108142.664
$ws.Range("J52").Value = 109771.8
$ws.Range("L52").Value = 109771.8
$ws.Range("N52").Value = -110359.8
$ws.Range("H62").Value = 3621
$ws.Range("I62").Value = 3899.6667
$ws.Range("J62").Value = 3481.6667
$ws.Range("K62").Value = 3899.6667
$ws.Range("L62").Value = 3481.6667
$ws.Range("M62").Value = -3275.6667
$ws.Range("N62").Value = -4729.6667
$ws.Range("H65").Value = 3621
$ws.Range("I65").Value = 3899.6667
$ws.Range("J65").Value = 3481.6667
$ws.Range("K65").Value = 19498.3335
$ws.Range("L65").Value = 17408.3335
$ws.Range("M65").Value = -16378.3335
$ws.Range("N65").Value = -23648.3335
$ws.Range("H99").Value = 9872.576999999999
$ws.Range("I99").Value = 9134.929
$ws.Range("J99").Value = 10733.167
$ws.Range("K99").Value = 9134.929
$ws.Range("L99").Value = 10733.167
$ws.Range("M99").Value = -7636.929
$ws.Range("N99").Value = -13729.167
$ws.Range("H126").Value = 9872.576999999999
$ws.Range("I126").Value = 9134.929
$ws.Range("J126").Value = 10733.167
$ws.Range("K126").Value = 27404.787
$ws.Range("L126").Value = 32199.501
$ws.Range("M126").Value = -24934.787
$ws.Range("N126").Value = -37139.501
$ws.Range("H132").Value = 3375.08
$ws.Range("I132").Value = 2469.6667
$ws.Range("K132").Value = 7409.000100000001
$ws.Range("M132").Value = -4879.000100000001
$ws.Range("H134").Value = 1045.8334
$ws.Range("I134").Value = 1045.8334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3137.5002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -602.5001999999999
$ws.Range("N134").Value = $null
$ws.Range("H141").Value = 466174.12
$ws.Range("J141").Value = 801995.2
$ws.Range("L141").Value = 801995.2
$ws.Range("N141").Value = -812355.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2693
$ws.Range("J7").Value = 370.2
$ws.Range("L7").Value = 1110.6
$ws.Range("N7").Value = -1334.6
$ws.Range("H15").Value = 295
$ws.Range("J15").Value = 500
$ws.Range("L15").Value = 1500
$ws.Range("N15").Value = -1780
$ws.Range("H92").Value = 506.33334
$ws.Range("I92").Value = 559.5
$ws.Range("K92").Value = 1678.5
$ws.Range("M92").Value = -430.5
$ws.Range("H131").Value = 39596.47
$ws.Range("I131").Value = 174497.17
$ws.Range("K131").Value = 523491.51
$ws.Range("M131").Value = -518451.51
$ws.Range("H137").Value = 8599.571
$ws.Range("J137").Value = 10499.6
$ws.Range("L137").Value = 31498.8
$ws.Range("N137").Value = -41698.8
$ws.Range("H139").Value = 4989.9375
$ws.Range("I139").Value = 3879.75
$ws.Range("J139").Value = 5360
$ws.Range("K139").Value = 11639.25
$ws.Range("L139").Value = 16080
$ws.Range("M139").Value = -6499.25
$ws.Range("N139").Value = -26360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 15765
$ws.Range("J47").Value = 15765
$ws.Range("L47").Value = 15765
$ws.Range("N47").Value = -16901
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 70
$ws.Range("N17").Value = $null
$ws.Range("H46").Value = 1114.4595
$ws.Range("I46").Value = 681.25
$ws.Range("J46").Value = 3887
$ws.Range("K46").Value = 681.25
$ws.Range("L46").Value = 3887
$ws.Range("M46").Value = -493.25
$ws.Range("N46").Value = -4263
$ws.Range("H61").Value = 5383.778
$ws.Range("I61").Value = 4590
$ws.Range("J61").Value = 6376
$ws.Range("K61").Value = 4590
$ws.Range("L61").Value = 6376
$ws.Range("M61").Value = -4388
$ws.Range("N61").Value = -6780
$ws.Range("H93").Value = 2443.889
$ws.Range("J93").Value = 4500
$ws.Range("L93").Value = 4500
$ws.Range("N93").Value = -6996
$ws.Range("H100").Value = 3998.8333
$ws.Range("I100").Value = 2998.5
$ws.Range("J100").Value = 5999.5
$ws.Range("K100").Value = 2998.5
$ws.Range("L100").Value = 5999.5
$ws.Range("M100").Value = -2457.5
$ws.Range("N100").Value = -7081.5
$ws.Range("H113").Value = 5383.778
$ws.Range("I113").Value = 4590
$ws.Range("J113").Value = 6376
$ws.Range("K113").Value = 4590
$ws.Range("L113").Value = 6376
$ws.Range("M113").Value = -2420
$ws.Range("N113").Value = -10716
$ws.Range("H132").Value = 90912216
$ws.Range("I132").Value = 2912.9333
$ws.Range("K132").Value = 8738.7999
$ws.Range("M132").Value = -6208.7999
$ws.Range("H136").Value = 2563.5
$ws.Range("I136").Value = 1766.3667
$ws.Range("K136").Value = 5299.1001
$ws.Range("M136").Value = -2749.1001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H46").Value = 93499.75
$ws.Range("J46").Value = 93499.75
$ws.Range("L46").Value = 93499.75
$ws.Range("N46").Value = -93961.75
$ws.Range("H55").Value = 8999.5
$ws.Range("I55").Value = 9999
$ws.Range("J55").Value = 8000
$ws.Range("K55").Value = 9999
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = -9722
$ws.Range("N55").Value = -8554
$ws.Range("H69").Value = 14999
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 14999
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 14999
$ws.Range("M69").Value = $null
$ws.Range("N69").Value = -16497
$ws.Range("H72").Value = 14999
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 14999
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 44997
$ws.Range("M72").Value = $null
$ws.Range("N72").Value = -52485
$ws.Range("H81").Value = 692.7692
$ws.Range("I81").Value = 692
$ws.Range("J81").Value = 697
$ws.Range("K81").Value = 1384
$ws.Range("L81").Value = 1394
$ws.Range("M81").Value = -323
$ws.Range("N81").Value = -3516
$ws.Range("H84").Value = 692.7692
$ws.Range("I84").Value = 692
$ws.Range("J84").Value = 697
$ws.Range("K84").Value = 6920
$ws.Range("L84").Value = 6970
$ws.Range("M84").Value = -1616
$ws.Range("N84").Value = -17578
$ws.Range("H112").Value = 26552.334
$ws.Range("J112").Value = 26552.334
$ws.Range("L112").Value = 26552.334
$ws.Range("N112").Value = -29506.334
$ws.Range("H132").Value = 3654
$ws.Range("I132").Value = 3725.3713
$ws.Range("J132").Value = 3237.6667
$ws.Range("K132").Value = 11176.1139
$ws.Range("L132").Value = 9713.000100000001
$ws.Range("M132").Value = -8646.1139
$ws.Range("N132").Value = -14773.0001
$ws.Range("H134").Value = 93499.75
$ws.Range("J134").Value = 93499.75
$ws.Range("L134").Value = 280499.25
$ws.Range("N134").Value = -285569.25
